$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 3.3.1 Security paragraph
# ---------------------------------------------------------------------------
Replace-Text "PHP’s PDO class extension. " "PHP’s internal libraries. "
Replace-Text "The exposure of passwords should be a concern" "The exposure of passwords is a concern"

# ---------------------------------------------------------------------------
# 3.3.2 Performance paragraph
# ---------------------------------------------------------------------------
Replace-Text "for data generation. Information should be displayed in real-time. The application will also ensure" "for data generation. The application shall also ensure"
Replace-Text "allow for quicker SQL query executions." "allow for quicker SQL query executions that are expected to return data within 100 nanoseconds. The performance shall also account for larger data execution, such as the schedule generator with a performance factor speed of no more than 5 seconds."

# ---------------------------------------------------------------------------
# 3.3.2 Cross-Browser Compatibility paragraph
# ---------------------------------------------------------------------------
Replace-Text "This application will ensure that the CSS " "This application shall ensure that the CSS "
Replace-Text " Ensuring browser compliance, the application will be able to run on" " Ensuring browser compliance, the application shall be able to run on"

# ---------------------------------------------------------------------------
# 3.3.3 Ease of Use paragraph
# ---------------------------------------------------------------------------
Replace-Text "The application will be straightforward" "The application shall be straightforward"
Replace-Text "their schedule immediately." "their schedule planner immediately."
Replace-Text "the data they are working on." "the data they were working on."

# ---------------------------------------------------------------------------
# First-line indents (0.5in / 720 twips = 36pt) on the four body paragraphs
# ---------------------------------------------------------------------------
$d.Paragraphs(5).Range.ParagraphFormat.FirstLineIndent = 36
$d.Paragraphs(9).Range.ParagraphFormat.FirstLineIndent = 36
$d.Paragraphs(13).Range.ParagraphFormat.FirstLineIndent = 36
$d.Paragraphs(17).Range.ParagraphFormat.FirstLineIndent = 36

# ---------------------------------------------------------------------------
# Move the _GoBack bookmark from its own empty paragraph to the start of the
# "The application shall be straightforward..." paragraph. Re-adding a
# bookmark with the same name relocates it (the old one is discarded), which
# also turns its former (now empty) host paragraph into a plain empty <w:p/>.
# ---------------------------------------------------------------------------
$easeOfUseRange = $d.Content.Duplicate
$easeOfUseRange.Find.Execute("The application shall be straightforward", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($easeOfUseRange.Start, $easeOfUseRange.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
